# Applies the ValueSet-split-method.xlsx metadata refresh:
#  - Version bumped 5.0.0 -> 6.0.0
#  - Date bumped to 2022-01-21T20:46:54+00:00
#  - Publisher value filled in ("Alvearie Team")
#  - Second duplicate "Contact" row removed, replaced by a single
#    "Jurisdiction" / "United States of America" row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" row (old row 11); this shifts every
# subsequent row up by one, turning the 15-row sheet into 14 rows.
$ws.Rows.Item(11).Delete()

# Version
$ws.Range("B3").Value = "6.0.0"

# Date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$ws.Range("B9").Value = "Alvearie Team"

# Former "Contact" row becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
